$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2 - habitat_type
$ws.Range("B2").Value2 = "402.213595695805" -as [double]
$ws.Range("D2").Value2 = "0.0000000000000000000000000000000000000000000000000000000000000000000000000000000000000924715937244439" -as [double]

# Row 3 - season
$ws.Range("B3").Value2 = "1097.0773183811" -as [double]
$ws.Range("D3").Value2 = "1.56726057069067e-237" -as [double]

# Row 4 - day_night
$ws.Range("B4").Value2 = "1663.51474739775" -as [double]

# Row 5 - habitat_type:season
$ws.Range("B5").Value2 = "592.590893772072" -as [double]
$ws.Range("D5").Value2 = "0.0000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000404885212428199" -as [double]

# Row 6 - habitat_type:day_night
$ws.Range("B6").Value2 = "152.891334737859" -as [double]
$ws.Range("D6").Value2 = "0.0000000000000000000000000146762662755666" -as [double]

# Row 7 - season:day_night
$ws.Range("B7").Value2 = "142.05602722908" -as [double]
$ws.Range("D7").Value2 = "0.0000000000000000000000000388067798537248" -as [double]

# Row 8 - habitat_type:season:day_night
$ws.Range("B8").Value2 = "108.781320209598" -as [double]
$ws.Range("D8").Value2 = "0.00000000171513043195422" -as [double]

$wb.Save()
